# Generate Report for Handback
# Refresh the timestamps recorded during the latest handback report
# generation run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview!G2 - "Latest HO Xliff Generate Date" for 38780d11-...md
$wsOverview.Range("G2").Value = "2016-08-23 19:08:41"

# de-de!H2 - "Correspond Handoff Datetime" for 38780d11-...md
# (shared the same original text/shared-string as Overview!G2, so it
# moves in lock-step with it)
$wsDeDe.Range("H2").Value = "2016-08-23 19:08:41"

# zh-cn!H2 - "Correspond Handoff Datetime" for 38780d11-...md
$wsZhCn.Range("H2").Value = "2016-08-23 19:08:36"

# zh-cn!K2 - "Correspond Handback DateTime" for 38780d11-...md
$wsZhCn.Range("K2").Value = "2016-08-23 19:09:06"

# de-de!K2 - "Correspond Handback DateTime" for 38780d11-...md
$wsDeDe.Range("K2").Value = "2016-08-23 19:09:16"
